$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(23, 3).Value = [double]"1.56845045540436e-265"
$ws.Cells.Item(24, 3).Value = [double]"7.050042827396416e-223"
$ws.Cells.Item(25, 3).Value = [double]"2.829392326046491e-192"
$ws.Cells.Item(26, 3).Value = [double]"3.180910882548996e-169"
$ws.Cells.Item(27, 3).Value = [double]"3.195165576662464e-151"
$ws.Cells.Item(28, 3).Value = [double]"9.18201336359421e-137"
$ws.Cells.Item(29, 3).Value = [double]"6.883522888340067e-125"
$ws.Cells.Item(30, 3).Value = [double]"5.892533096934832e-115"
$ws.Cells.Item(31, 3).Value = [double]"1.603772713037502e-106"
$ws.Cells.Item(32, 3).Value = [double]"2.887841327313808e-99"
$ws.Cells.Item(33, 3).Value = [double]"5.893601627806694e-93"
$ws.Cells.Item(34, 3).Value = [double]"2.042721823947946e-87"
$ws.Cells.Item(35, 3).Value = [double]"3.876186715972135e-78"
$ws.Cells.Item(36, 3).Value = [double]"1.146494151266776e-70"
$ws.Cells.Item(37, 3).Value = [double]"4.479997314112631e-57"
$ws.Cells.Item(38, 3).Value = [double]"7.213750392475025e-48"
$ws.Cells.Item(39, 3).Value = [double]"3.464307538328817e-41"
$ws.Cells.Item(40, 3).Value = [double]"4.278630249216903e-36"
$ws.Cells.Item(41, 3).Value = [double]"4.590752145370398e-32"
$ws.Cells.Item(42, 3).Value = [double]"9.066317752488081e-29"
$ws.Cells.Item(43, 3).Value = [double]"1.258585094456025e-23"
$ws.Cells.Item(44, 3).Value = [double]"9.426703589427085e-20"
$ws.Cells.Item(45, 3).Value = [double]"1.018318890387792e-16"
$ws.Cells.Item(46, 3).Value = [double]"2.785879063589272e-14"
$ws.Cells.Item(47, 3).Value = [double]"2.743822352787366e-12"
$ws.Cells.Item(48, 3).Value = [double]"1.253706400610929e-08"
$ws.Cells.Item(49, 3).Value = [double]"3.718615117761819e-06"
$ws.Cells.Item(50, 3).Value = [double]"0.0002213651962603142"
$ws.Cells.Item(51, 3).Value = [double]"0.004771456859910801"
$ws.Cells.Item(52, 3).Value = [double]"0.3503997549335854"
$ws.Cells.Item(53, 3).Value = [double]"6.063032433126999"
$ws.Cells.Item(54, 3).Value = [double]"45.70338917419514"
$ws.Cells.Item(55, 3).Value = [double]"204.6695357129844"
$ws.Cells.Item(56, 3).Value = [double]"1607.720445086961"
$ws.Cells.Item(57, 3).Value = [double]"6108.157404421996"
$ws.Cells.Item(58, 3).Value = [double]"15374.56844557216"
$ws.Cells.Item(59, 3).Value = [double]"30001.25528387708"
$ws.Cells.Item(60, 3).Value = [double]"49509.79880754172"
$ws.Cells.Item(61, 3).Value = [double]"72776.89912643001"
